$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("15:17").Insert()
Write-Host "AutoFilterMode after insert:" $ws.AutoFilterMode
if ($ws.AutoFilter) {
    Write-Host "range still:" 
}
